$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in T5 (F) and T6 (G) grades for each student row (2-6)
$ws.Range("F2").Value = 1.25
$ws.Range("G2").Value = 1

$ws.Range("F3").Value = 1.25
$ws.Range("G3").Value = 1.25

$ws.Range("F4").Value = 1.25
$ws.Range("G4").Value = 1

$ws.Range("F5").Value = 1.25
$ws.Range("G5").Value = 1

$ws.Range("F6").Value = 1.25
$ws.Range("G6").Value = 1.25

# Update selection to match the author's final cursor position
$ws.Range("G7").Select()
